# "añado reseña nuestra aplicacion"
# Fill in the SUS questionnaire answers (columns C and E) for the 10
# questions in rows 14-23. Column D already holds a formula that derives
# from column C, so it (and the row-24 totals) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$answers = @{
    14 = @{ C = 2; E = 2 }
    15 = @{ C = 1; E = 1 }
    16 = @{ C = 3; E = 4 }
    17 = @{ C = 2; E = 1 }
    18 = @{ C = 4; E = 4 }
    19 = @{ C = 2; E = 2 }
    20 = @{ C = 3; E = 4 }
    21 = @{ C = 1; E = 1 }
    22 = @{ C = 2; E = 4 }
    23 = @{ C = 1; E = 1 }
}

foreach ($row in $answers.Keys | Sort-Object) {
    $ws.Range("C$row").Value = $answers[$row].C
    $ws.Range("E$row").Value = $answers[$row].E
}

# Scroll the view up a bit and leave the selection on the new "Normalizado"
# total for the second reviewer, matching the reviewer's final cursor spot.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E24").Select()
